# Update countries & provincias Spain
# - Refresh the COVID-19 numbers for several countries.
# - Dinamarca overtakes Noruega, and Filipinas overtakes Mexico / Arabia
#   Saudita in the ranking, so those rows swap country labels while the
#   row below picks up the displaced country's previous figures.
# - Refresh the "Datos actualizados" timestamp cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: Polonia (values refreshed, same country) ---
$ws.Range("B30").Value = 7408
$ws.Range("C30").Value = 206
$ws.Range("E30").Value = 6472
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 268

# --- Row 32: now Dinamarca (was Noruega) ---
$ws.Range("A32").Value = "Dinamarca"
$ws.Range("B32").Value = 6681
$ws.Range("C32").Value = 170
$ws.Range("D32").Value = 2515
$ws.Range("E32").Value = 3867
$ws.Range("F32").Value = 100
$ws.Range("H32").Value = 299

# --- Row 33: now Noruega (was Dinamarca) ---
$ws.Range("A33").Value = "Noruega"
$ws.Range("B33").Value = 6623
$ws.Range("D33").Value = 32
$ws.Range("E33").Value = 6452
$ws.Range("F33").Value = 59
$ws.Range("H33").Value = 139

# --- Row 37: now Filipinas (was Mexico) ---
$ws.Range("A37").Value = "Filipinas"
$ws.Range("B37").Value = 5453
$ws.Range("C37").Value = 230
$ws.Range("D37").Value = 353
$ws.Range("E37").Value = 4751
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 14
$ws.Range("H37").Value = 349

# --- Row 38: now Mexico (was Arabia Saudita) ---
$ws.Range("A38").Value = "Mexico"
$ws.Range("B38").Value = 5399
$ws.Range("C38").Value = 385
$ws.Range("D38").Value = 2125
$ws.Range("E38").Value = 2868
$ws.Range("F38").Value = 207
$ws.Range("G38").Value = 74
$ws.Range("H38").Value = 406

# --- Row 39: now Arabia Saudita (was Filipinas) ---
$ws.Range("A39").Value = "Arabia Saudita"
$ws.Range("B39").Value = 5369
$ws.Range("D39").Value = 889
$ws.Range("E39").Value = 4407
$ws.Range("F39").Value = 59
$ws.Range("H39").Value = 73

# --- Row 90: Libano (values refreshed, same country) ---
$ws.Range("D90").Value = 81
$ws.Range("E90").Value = 539

# --- A1: refresh the "updated at" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 10:22"
